$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("linea")

# Selection change (best-effort, may not affect saved XML depending on runtime)
$ws.Range("L18").Select()

# K column (CAPATAZ) updates
$ws.Range("K5").Value = "DAVID"
$ws.Range("K11").Value = "PEPE"
$ws.Range("K12").Value = "PEPE"

# L column (OBSERVACIONES) updates
$ws.Range("L5").Value = "UNO"
$ws.Range("L6").Value = "DOS"
$ws.Range("L7").Value = "TRES"
$ws.Range("L8").Value = "CUATRO"
$ws.Range("L9").Value = "CINCO"
$ws.Range("L10").Value = "SEIS"
$ws.Range("L11").Value = "SIETE"
$ws.Range("L12").Value = "OCHO"
$ws.Range("L13").Value = "NUEVE"
$ws.Range("L14").Value = "DIEZ"
$ws.Range("L15").Value = "ONCE"
$ws.Range("L16").Value = "DOCE"
$ws.Range("L17").Value = "TRECE"
$ws.Range("L18").Value = "CATORCE"

$wb.Save()
